{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the \"LOB1012: ...\" requirement paragraph, then delete the three\n// paragraphs that used to follow it (the blank spacer line, the\n// \"Ver no Jupiter ...\" line, and the \"\u00a9 2020 ...\" copyright line) which made\n// up the old page-footer block scraped from the site.\nconst items = paragraphs.items;\nlet targetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"LOB1012: Estat\u00edstica (Requisito)\") {\n    targetIndex = i;\n    break;\n  }\n}\n\nif (targetIndex !== -1) {\n  // Delete in reverse order so earlier indices stay valid.\n  for (let i = targetIndex + 3; i >= targetIndex + 1; i--) {\n    items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph that contains the exact \"LOB1012: ...\" requirement line,\n# then delete the three paragraphs that followed it in the old \"Jupiter footer\"\n# block: the blank spacer paragraph, the \"Ver no Jupiter ...\" line, and the\n# \"\u00a9 2020 ...\" copyright line. Walking by content (rather than a hard-coded\n# index) keeps this robust if earlier paragraphs in the document shift.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"LOB1012: Estat\u00edstica (Requisito)\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # The paragraph immediately after the requisito line is the blank spacer;\n    # its two successors are the \"Ver no Jupiter...\" and \"\u00a9 2020...\" lines.\n    $blank = $target.Next()\n    $jupiter = $blank.Next()\n    $copyright = $jupiter.Next()\n\n    $copyright.Range.Delete()\n    $jupiter.Range.Delete()\n    $blank.Range.Delete()\n}\n"}
